# Add two new columns I (I0) and J (IF) to the worksheet, mirroring the
# formatting of the existing header/data columns (e.g. column H / "IP").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new header cells I1 and J1 ---------------------
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the existing header formatting (bold font, borders, centered
# alignment) from H1, which already carries the header style, onto the
# two new header cells so they match the rest of the header row.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# --- Data rows (rows 2-12): new numeric values for I and J --------------
$values = @{
    2  = @(1, 3)
    3  = @(1, 4)
    4  = @(1, 4)
    5  = @(1, 5)
    6  = @(1, 5)
    7  = @(1, 5)
    8  = @(5, 9)
    9  = @(1, 4)
    10 = @(1, 5)
    11 = @(8, 8)
    12 = @(1, 2)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]  # column I
    $ws.Cells.Item($row, 10).Value = $pair[1] # column J
}
